$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 493.1
$ws.Range("I41").Value = 579.8570999999999
$ws.Range("K41").Value = 579.8570999999999
$ws.Range("M41").Value = -139.8570999999999

$ws.Range("H64").Value = 5016
$ws.Range("I64").Value = 4984
$ws.Range("K64").Value = 4984
$ws.Range("M64").Value = -4736

$ws.Range("H67").Value = 5016
$ws.Range("I67").Value = 4984
$ws.Range("K67").Value = 4984
$ws.Range("M67").Value = -4126

$ws.Range("H115").Value = 1105.6364
$ws.Range("I115").Value = 424.25
$ws.Range("J115").Value = 1495
$ws.Range("K115").Value = 1272.75
$ws.Range("L115").Value = 4485
$ws.Range("M115").Value = 294.25
$ws.Range("N115").Value = -7619

$ws.Range("H129").Value = 902.4167
$ws.Range("I129").Value = 902.4167
$ws.Range("K129").Value = 2707.2501
$ws.Range("M129").Value = 2292.7499

$ws.Range("H131").Value = 3018.7778
$ws.Range("I131").Value = 896.125
$ws.Range("K131").Value = 2688.375
$ws.Range("M131").Value = 2351.625

$ws.Range("H135").Value = 1648.7587
$ws.Range("J135").Value = 2775.7
$ws.Range("L135").Value = 24981.3
$ws.Range("N135").Value = -30051.3

$ws.Range("H137").Value = 1737.421
$ws.Range("I137").Value = 976.1539
$ws.Range("J137").Value = 3386.8333
$ws.Range("K137").Value = 2928.4617
$ws.Range("L137").Value = 10160.4999
$ws.Range("M137").Value = -378.4616999999998
$ws.Range("N137").Value = -15260.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644.96155
$ws.Range("I2").Value = 545
$ws.Range("J2").Value = 1194.75
$ws.Range("K2").Value = 545
$ws.Range("L2").Value = 1194.75
$ws.Range("M2").Value = -432
$ws.Range("N2").Value = -1420.75

$ws.Range("H32").Value = 40322.242
$ws.Range("I32").Value = 44134.24
$ws.Range("J32").Value = 16497.25
$ws.Range("K32").Value = 44134.24
$ws.Range("L32").Value = 16497.25
$ws.Range("M32").Value = -43847.24
$ws.Range("N32").Value = -17071.25

$ws.Range("H74").Value = 35664.45
$ws.Range("I74").Value = 37972.926
$ws.Range("K74").Value = 37972.926
$ws.Range("M74").Value = -37098.926

$ws.Range("H77").Value = 35664.45
$ws.Range("I77").Value = 37972.926
$ws.Range("K77").Value = 189864.63
$ws.Range("M77").Value = -185496.63

$ws.Range("H102").Value = 3113.9048
$ws.Range("I102").Value = 2673.8125
$ws.Range("K102").Value = 2673.8125
$ws.Range("M102").Value = -1051.8125

$ws.Range("H116").Value = 644.96155
$ws.Range("I116").Value = 545
$ws.Range("J116").Value = 1194.75
$ws.Range("K116").Value = 545
$ws.Range("L116").Value = 1194.75
$ws.Range("M116").Value = 1749
$ws.Range("N116").Value = -5782.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644.96155
$ws.Range("I3").Value = 545
$ws.Range("J3").Value = 1194.75
$ws.Range("K3").Value = 545
$ws.Range("L3").Value = 1194.75
$ws.Range("M3").Value = -431
$ws.Range("N3").Value = -1422.75

$ws.Range("H80").Value = 469.2
$ws.Range("J80").Value = 360.75
$ws.Range("L80").Value = 360.75
$ws.Range("N80").Value = -2356.75

$ws.Range("H83").Value = 469.2
$ws.Range("J83").Value = 360.75
$ws.Range("L83").Value = 1803.75
$ws.Range("N83").Value = -11787.75

$ws.Range("H105").Value = 3142.261
$ws.Range("I105").Value = 3337.4443
$ws.Range("J105").Value = 2439.6
$ws.Range("K105").Value = 3337.4443
$ws.Range("L105").Value = 2439.6
$ws.Range("M105").Value = -1590.4443
$ws.Range("N105").Value = -5933.6

$ws.Range("H107").Value = 1149.6
$ws.Range("I107").Value = 1149.6
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1149.6
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 770.4000000000001
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 1691.7778
$ws.Range("I134").Value = 1703.6471
$ws.Range("K134").Value = 5110.9413
$ws.Range("M134").Value = -2575.9413

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 545.4
$ws.Range("I22").Value = 342.9
$ws.Range("K22").Value = 342.9
$ws.Range("M22").Value = 7.100000000000023

$ws.Range("H31").Value = 2537.9546
$ws.Range("I31").Value = 1791.5385
$ws.Range("J31").Value = 3616.111
$ws.Range("K31").Value = 1791.5385
$ws.Range("L31").Value = 3616.111
$ws.Range("M31").Value = -1496.5385
$ws.Range("N31").Value = -4206.111

$ws.Range("H34").Value = 2537.9546
$ws.Range("I34").Value = 1791.5385
$ws.Range("J34").Value = 3616.111
$ws.Range("K34").Value = 1791.5385
$ws.Range("L34").Value = 3616.111
$ws.Range("M34").Value = -1589.5385
$ws.Range("N34").Value = -4020.111

$ws.Range("H58").Value = 102429.7
$ws.Range("I58").Value = 127160.375
$ws.Range("J58").Value = 3507
$ws.Range("K58").Value = 127160.375
$ws.Range("L58").Value = 3507
$ws.Range("M58").Value = -126957.375
$ws.Range("N58").Value = -3913

$ws.Range("H107").Value = 464.875
$ws.Range("I107").Value = 452.53333
$ws.Range("J107").Value = 650
$ws.Range("K107").Value = 452.53333
$ws.Range("L107").Value = 650
$ws.Range("M107").Value = 1467.46667
$ws.Range("N107").Value = -4490

$ws.Range("H132").Value = 999.89746
$ws.Range("I132").Value = 1035.7576
$ws.Range("J132").Value = 802.6667
$ws.Range("K132").Value = 3107.2728
$ws.Range("L132").Value = 2408.0001
$ws.Range("M132").Value = -577.2727999999997
$ws.Range("N132").Value = -7468.0001

$ws.Range("H134").Value = 40527.58
$ws.Range("I134").Value = 50993.75
$ws.Range("J134").Value = 5640.3335
$ws.Range("K134").Value = 152981.25
$ws.Range("L134").Value = 16921.0005
$ws.Range("M134").Value = -150446.25
$ws.Range("N134").Value = -21991.0005

$ws.Range("H136").Value = 102429.7
$ws.Range("I136").Value = 127160.375
$ws.Range("J136").Value = 3507
$ws.Range("K136").Value = 381481.125
$ws.Range("L136").Value = 10521
$ws.Range("M136").Value = -378931.125
$ws.Range("N136").Value = -15621

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 960.4286
$ws.Range("I113").Value = 1065
$ws.Range("J113").Value = 902.3333
$ws.Range("K113").Value = 3195
$ws.Range("L113").Value = 2706.9999
$ws.Range("M113").Value = -1025
$ws.Range("N113").Value = -7046.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 44393.418
$ws.Range("I132").Value = 61293.94
$ws.Range("K132").Value = 183881.82
$ws.Range("M132").Value = -181351.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 17522.5
$ws.Range("I45").Value = 24999
$ws.Range("J45").Value = 10046
$ws.Range("K45").Value = 24999
$ws.Range("L45").Value = 10046
$ws.Range("M45").Value = -24592
$ws.Range("N45").Value = -10860

$ws.Range("H46").Value = 53231.668
$ws.Range("I46").Value = 63078
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 63078
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -62890
$ws.Range("N46").Value = -4376

$ws.Range("H55").Value = 1232.6666
$ws.Range("I55").Value = 879.9286
$ws.Range("J55").Value = 1938.1428
$ws.Range("K55").Value = 879.9286
$ws.Range("L55").Value = 1938.1428
$ws.Range("M55").Value = -706.9286
$ws.Range("N55").Value = -2284.1428

$ws.Range("H122").Value = 4216.4
$ws.Range("I122").Value = 3050
$ws.Range("J122").Value = 4994
$ws.Range("K122").Value = 9150
$ws.Range("L122").Value = 14982
$ws.Range("M122").Value = -6700
$ws.Range("N122").Value = -19882

$ws.Range("H136").Value = 2243.5
$ws.Range("I136").Value = 2243.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6730.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4180.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 53024.184
$ws.Range("I126").Value = 70146.125
$ws.Range("K126").Value = 210438.375
$ws.Range("M126").Value = -207968.375

$ws.Range("H132").Value = 27289.605
$ws.Range("I132").Value = 28289.451
$ws.Range("K132").Value = 84868.353
$ws.Range("M132").Value = -82338.353

$ws.Range("H136").Value = 3409.4443
$ws.Range("I136").Value = 3340
$ws.Range("J136").Value = 3496.25
$ws.Range("K136").Value = 10020
$ws.Range("L136").Value = 10488.75
$ws.Range("M136").Value = -7470
$ws.Range("N136").Value = -15588.75
